$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.250.90'
$ws.Range('E2').Value = '  -0.54%  '

$ws.Range('D3').Value = '1.785.57'
$ws.Range('E3').Value = '  -1.71%  '

$ws.Range('E4').Value = '  +0.16%  '

$ws.Range('D5').Value = "'334.39"
$ws.Range('E5').Value = '  -2.96%  '

$ws.Range('E6').Value = '  +0.10%  '

$ws.Range('D7').Value = "'0.3789"
$ws.Range('E7').Value = '  -1.60%  '

$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').Value = "'48.78"
$ws.Range('E8').Value = '  -3.26%  '

$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').Value = "'0.3430"
$ws.Range('E9').Value = '  -2.97%  '

$ws.Range('D10').Value = "'1.198"
$ws.Range('E10').Value = '  -3.57%  '

$ws.Range('D11').Value = "'0.07494"
$ws.Range('E11').Value = '  -3.63%  '

$ws.Range('E12').Value = '  +0.02%  '

$ws.Range('D13').Value = "'21.90"
$ws.Range('E13').Value = '  -3.43%  '

$ws.Range('D14').Value = "'6.475"
$ws.Range('E14').Value = '  -2.62%  '

$ws.Range('D15').Value = '1.783.25'
$ws.Range('E15').Value = '  -1.89%  '

$ws.Range('D16').Value = "'7.089"
$ws.Range('E16').Value = '  -2.29%  '

$ws.Range('E17').Value = '  -2.65%  '

$ws.Range('D18').Value = "'0.06650"
$ws.Range('E18').Value = '  -2.16%  '

$ws.Range('D19').Value = "'83.85"
$ws.Range('E19').Value = '  -3.97%  '

$ws.Range('E20').Value = '  +0.14%  '

$ws.Range('D21').Value = "'6.642"
$ws.Range('E21').Value = '  +1.10%  '

$ws.Range('D22').Value = "'17.36"
$ws.Range('E22').Value = '  -3.09%  '

$ws.Range('D23').Value = '27.255.59'
$ws.Range('E23').Value = '  -0.49%  '

$ws.Range('D24').Value = "'12.37"
$ws.Range('E24').Value = '  -6.21%  '

$ws.Range('D25').Value = "'2.412"
$ws.Range('E25').Value = '  -2.55%  '

$ws.Range('D26').Value = "'1.509"
$ws.Range('E26').Value = '  +0.32%  '

$ws.Range('D27').Value = "'2.544"
$ws.Range('E27').Value = '  -7.03%  '

$ws.Range('D28').Value = "'21.35"
$ws.Range('E28').Value = '  -3.80%  '

$ws.Range('D29').Value = "'152.96"
$ws.Range('E29').Value = '  -1.11%  '

$ws.Range('D30').Value = '1.988.89'
$ws.Range('E30').Value = '  -1.59%  '

$ws.Range('D31').Value = "'134.09"
$ws.Range('E31').Value = '  -2.30%  '

$ws.Range('D32').Value = "'4.032"
$ws.Range('E32').Value = '  -2.19%  '

$ws.Range('D33').Value = "'6.089"
$ws.Range('E33').Value = '  -5.28%  '

$ws.Range('D34').Value = "'0.08700"
$ws.Range('E34').Value = '  -1.67%  '

$ws.Range('D35').Value = "'13.24"
$ws.Range('E35').Value = '  -4.66%  '

$ws.Range('D36').Value = "'1.669"
$ws.Range('E36').Value = '  -3.11%  '

$ws.Range('D37').Value = "'0.6945"
$ws.Range('E37').Value = '  -2.35%  '

$ws.Range('D38').Value = "'5.451"
$ws.Range('E38').Value = '  -4.01%  '

$ws.Range('D39').Value = "'0.2207"
$ws.Range('E39').Value = '  -2.99%  '

$ws.Range('D40').Value = "'8.823"
$ws.Range('E40').Value = '  -2.68%  '

$ws.Range('D41').Value = "'0.06331"
$ws.Range('E41').Value = '  -3.93%  '

$ws.Range('D42').Value = "'0.02337"
$ws.Range('E42').Value = '  -3.53%  '

$ws.Range('E43').Value = '  -1.77%  '

$ws.Range('D44').Value = "'14.47"
$ws.Range('E44').Value = '  -3.96%  '

$ws.Range('D45').Value = "'0.6534"
$ws.Range('E45').Value = '  -1.75%  '

$ws.Range('D46').Value = "'1.001"
$ws.Range('E46').Value = '  +0.09%  '

$ws.Range('D47').Value = "'3.840"
$ws.Range('E47').Value = '  -3.52%  '

$ws.Range('D48').Value = "'2.152"
$ws.Range('E48').Value = '  -2.18%  '

$ws.Range('D49').Value = "'129.36"
$ws.Range('E49').Value = '  -2.94%  '

$ws.Range('D50').Value = "'0.07133"
$ws.Range('E50').Value = '  -3.32%  '

$ws.Range('D51').Value = "'79.17"
$ws.Range('E51').Value = '  -2.40%  '
